$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 98; $r++) {
    $ws.Cells.Item($r, 11).Formula = "=E$r/D$r"
    $ws.Cells.Item($r, 12).Formula = "=H$r/F$r"
}

$ws.Range("K1:L98").Select()
$excel.ActiveWindow.ScrollColumn = 2
